$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A45").Value = "29-10-2025"
$ws.Range("B45").Value = "The price of gold in India today is ₹12,158 per gram for 24 karat gold, ₹11,145 per gram for 22 karat gold and ₹9,119 per gram for 18 karat gold (also called 999 gold)."
